# Update the cryptos price list (Coin / Link / Price / Volume(1h) columns).
# All data cells in this sheet are stored as plain text, so for any "Price"
# value that looks like a bare number (e.g. "1.00", "233.05") we prefix it
# with a leading apostrophe before assigning it. That forces Excel to keep
# the value as text (preserving trailing zeros / exact digits) instead of
# silently re-interpreting it as a numeric value, which would both change
# the cell type and risk floating point rounding artifacts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    if ($value -match '^-?\d+(\.\d+)?$') {
        $range.Value = "'" + $value
    } else {
        $range.Value = $value
    }
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "41.934.70"
$ws.Range("E2").Value = "  +1.39%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "2.233.55"
$ws.Range("E3").Value = "  -0.27%  "

# Row 4 - TetherUSD
Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  -0.39%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "233.05"
$ws.Range("E5").Value = "  +1.90%  "

# Row 6 - XRP
Set-TextValue $ws.Range("D6") "0.622"
$ws.Range("E6").Value = "  -2.25%  "

# Row 7 - Solana
Set-TextValue $ws.Range("D7") "60.78"
$ws.Range("E7").Value = "  -6.38%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.04%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -0.32%  "

# Row 10 - OKB
Set-TextValue $ws.Range("D10") "58.19"
$ws.Range("E10").Value = "  -4.28%  "

# Row 11 - Dogecoin
Set-TextValue $ws.Range("D11") "0.0901"
$ws.Range("E11").Value = "  +1.82%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  -0.58%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D13") "2.565.57"
$ws.Range("E13").Value = "  -0.54%  "

# Row 14 - Chainlink
Set-TextValue $ws.Range("D14") "15.54"
$ws.Range("E14").Value = "  -3.50%  "

# Row 15 - Avalanche
$ws.Range("E15").Value = "  +1.82%  "

# Row 16 - Polkadot
Set-TextValue $ws.Range("D16") "5.68"
$ws.Range("E16").Value = "  +0.71%  "

# Row 17 - Polygon
Set-TextValue $ws.Range("D17") "0.804"
$ws.Range("E17").Value = "  -2.94%  "

# Row 18 - WrappedEther
Set-TextValue $ws.Range("D18") "2.243.56"
$ws.Range("E18").Value = "  -0.07%  "

# Row 19 - WrappedBTC
Set-TextValue $ws.Range("D19") "41.808.79"
$ws.Range("E19").Value = "  +1.39%  "

# Row 20 - ShibaInu
Set-TextValue $ws.Range("D20") "0.0₃0911"
$ws.Range("E20").Value = "  +0.23%  "

# Row 21 - Litecoin
Set-TextValue $ws.Range("D21") "72.67"
$ws.Range("E21").Value = "  -1.71%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +0.70%  "

# Row 23 - BitcoinCash
Set-TextValue $ws.Range("D23") "248.52"
$ws.Range("E23").Value = "  -2.63%  "

# Row 24 - Dai
Set-TextValue $ws.Range("D24") "1.00"
$ws.Range("E24").Value = "  -0.06%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  +0.25%  "

# Row 26 - Toncoin
Set-TextValue $ws.Range("D26") "2.31"
$ws.Range("E26").Value = "  -3.57%  "

# Row 27 - Cosmos
Set-TextValue $ws.Range("D27") "9.72"
$ws.Range("E27").Value = "  -0.45%  "

# Row 28 - Monero
Set-TextValue $ws.Range("D28") "169.67"
$ws.Range("E28").Value = "  -1.93%  "

# Row 29 - Kaspa
$ws.Range("E29").Value = "  -2.54%  "

# Row 30 - EthereumClassic
Set-TextValue $ws.Range("D30") "19.95"
$ws.Range("E30").Value = "  -2.42%  "

# Row 31 - ImmutableX
$ws.Range("E31").Value = "  -2.69%  "

# Row 32 - WEMIXToken
$ws.Range("E32").Value = "  -9.50%  "

# Row 33 - Stellar
$ws.Range("E33").Value = "  -1.59%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  +3.48%  "

# Row 35 - Filecoin
$ws.Range("E35").Value = "  +0.72%  "

# Row 36 - Hedera
$ws.Range("E36").Value = "  +4.43%  "

# Row 37 - THORChain
Set-TextValue $ws.Range("D37") "6.60"
$ws.Range("E37").Value = "  -8.51%  "

# Row 38 - LidoDAOToken
$ws.Range("E38").Value = "  -2.28%  "

# Row 39 - RenderToken
$ws.Range("E39").Value = "  -5.37%  "

# Row 40 - TerraClassic
Set-TextValue $ws.Range("D40") "0.000246"
$ws.Range("E40").Value = "  +3.77%  "

# Row 41 - BinanceUSD
Set-TextValue $ws.Range("D41") "1.00"
$ws.Range("E41").Value = "  -0.11%  "

# Row 42 - VeChain
$ws.Range("E42").Value = "  +1.92%  "

# Row 43 - FraxShare
Set-TextValue $ws.Range("D43") "8.67"
$ws.Range("E43").Value = "  -0.94%  "

# Row 44 - TrustWalletToken
$ws.Range("E44").Value = "  -0.39%  "

# Row 45 - now Aave (was FTXToken; rows 45/46 swapped order)
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D45") "99.13"
$ws.Range("E45").Value = "  -2.60%  "

# Row 46 - now FTXToken (was Aave)
$ws.Range("B46").Value = "FTXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws.Range("D46") "4.49"
$ws.Range("E46").Value = "  -7.73%  "

# Row 47 - Cronos
Set-TextValue $ws.Range("D47") "0.0966"
$ws.Range("E47").Value = "  +2.65%  "

# Row 48 - Maker
Set-TextValue $ws.Range("D48") "1.473.86"
$ws.Range("E48").Value = "  -2.54%  "

# Row 49 - InjectiveProtocol
Set-TextValue $ws.Range("D49") "16.67"
$ws.Range("E49").Value = "  -6.34%  "

# Row 50 - now NEARProtocol (was HuobiToken; rows 50/51 swapped order)
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D50") "2.28"
$ws.Range("E50").Value = "  +8.02%  "

# Row 51 - now HuobiToken (was NEARProtocol)
$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D51") "2.77"
$ws.Range("E51").Value = "  -2.31%  "
